$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Unprotect sheet (it's protected with password hash D382) so values can be edited
$ws.Unprotect()

# Update the confidential notice date string in A9
$ws.Range("A9").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-05-05 for illustrative purposes only and are subject to change."

# Update Weight (D) and Percent Change (E) values for rows 2-6
$ws.Range("D2").Value = 0.2533359143840864
$ws.Range("E2").Value = 0.01267605633802815

$ws.Range("D3").Value = 0.2520534700962468
$ws.Range("E3").Value = 0.008167710318540689

$ws.Range("D4").Value = 0.2465592245424962
$ws.Range("E4").Value = 0.0009668374746203412

$ws.Range("D5").Value = 0.2480513909771706
$ws.Range("E5").Value = -0.001949571094359337

$ws.Range("E6").Value = 0.005024788927944224

# Re-protect the sheet (restores protected state; exact legacy password hash
# cannot be reproduced via COM, but protection itself is restored)
$ws.Protect()
